$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44188
$ws.Range("M2").Value = 30

# Row 3
$ws.Range("D3").Value = 44175
$ws.Range("M3").Value = 25
$ws.Range("N3").Value = 20000
$ws.Range("O3").Value = 20000
$ws.Range("P3").Value = 20000
$ws.Range("S3").Value = 4000

# Row 5
$ws.Range("D5").Value = 44193
$ws.Range("M5").Value = 40
$ws.Range("N5").Value = 15000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 15000
$ws.Range("S5").Value = 3000

# Row 8
$ws.Range("D8").Value = 44189

# Row 9
$ws.Range("D9").Value = 44186
